$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list) - update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 95
$wsExhibit.Range("F4").Value = 1062

# Sheet "全部类型" (all types) - same rows mirrored here, update likewise
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 95
$wsAll.Range("F4").Value = 1062
